# Insert 6 new weekly price rows at the top of the Chirimoya data block
# (rows 328-333), pushing the existing rows down by 6 (328->334 ... 348->354).
# Mirrors a new week's report ("Fruta / hortaliza, semanal") being
# prepended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert six blank rows before the current row 328; Excel's InsertRows
# pushes everything from the old row 328 onward down to row 334 onward,
# and copies the row-328 formatting (incl. the date style on column D)
# down into the freshly inserted rows.
$ws.Rows("328:333").Insert()

# New row 328: Especial, Cultivar IV Región, Provincia del Elquí
$ws.Cells.Item(328, 1).Value = 3
$ws.Cells.Item(328, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(328, 3).Value = "Coquimbo"
$ws.Cells.Item(328, 4).Value = 45166
$ws.Cells.Item(328, 5).Value = 5
$ws.Cells.Item(328, 6).Value = "Fruta"
$ws.Cells.Item(328, 7).Value = 100107
$ws.Cells.Item(328, 8).Value = "Otros"
$ws.Cells.Item(328, 9).Value = 100107002
$ws.Cells.Item(328, 10).Value = "Chirimoya"
$ws.Cells.Item(328, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(328, 12).Value = "Especial"
$ws.Cells.Item(328, 13).Value = 65
$ws.Cells.Item(328, 14).Value = 32000
$ws.Cells.Item(328, 15).Value = 32000
$ws.Cells.Item(328, 16).Value = 32000
$ws.Cells.Item(328, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(328, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(328, 19).Value = 3200
$ws.Cells.Item(328, 20).Value = 10

# New row 329: Primera, Cultivar IV Región, Provincia del Elquí
$ws.Cells.Item(329, 1).Value = 3
$ws.Cells.Item(329, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 45166
$ws.Cells.Item(329, 5).Value = 5
$ws.Cells.Item(329, 6).Value = "Fruta"
$ws.Cells.Item(329, 7).Value = 100107
$ws.Cells.Item(329, 8).Value = "Otros"
$ws.Cells.Item(329, 9).Value = 100107002
$ws.Cells.Item(329, 10).Value = "Chirimoya"
$ws.Cells.Item(329, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(329, 12).Value = "Primera"
$ws.Cells.Item(329, 13).Value = 68
$ws.Cells.Item(329, 14).Value = 30000
$ws.Cells.Item(329, 15).Value = 30000
$ws.Cells.Item(329, 16).Value = 30000
$ws.Cells.Item(329, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(329, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(329, 19).Value = 3000
$ws.Cells.Item(329, 20).Value = 10

# New row 330: Segunda, Cultivar IV Región, Provincia del Elquí
$ws.Cells.Item(330, 1).Value = 3
$ws.Cells.Item(330, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(330, 3).Value = "Coquimbo"
$ws.Cells.Item(330, 4).Value = 45166
$ws.Cells.Item(330, 5).Value = 5
$ws.Cells.Item(330, 6).Value = "Fruta"
$ws.Cells.Item(330, 7).Value = 100107
$ws.Cells.Item(330, 8).Value = "Otros"
$ws.Cells.Item(330, 9).Value = 100107002
$ws.Cells.Item(330, 10).Value = "Chirimoya"
$ws.Cells.Item(330, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(330, 12).Value = "Segunda"
$ws.Cells.Item(330, 13).Value = 57
$ws.Cells.Item(330, 14).Value = 28000
$ws.Cells.Item(330, 15).Value = 28000
$ws.Cells.Item(330, 16).Value = 28000
$ws.Cells.Item(330, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(330, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(330, 19).Value = 2800
$ws.Cells.Item(330, 20).Value = 10

# New row 331: Especial, Cultivar V Región, Provincia de Quillota
$ws.Cells.Item(331, 1).Value = 3
$ws.Cells.Item(331, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 45166
$ws.Cells.Item(331, 5).Value = 5
$ws.Cells.Item(331, 6).Value = "Fruta"
$ws.Cells.Item(331, 7).Value = 100107
$ws.Cells.Item(331, 8).Value = "Otros"
$ws.Cells.Item(331, 9).Value = 100107002
$ws.Cells.Item(331, 10).Value = "Chirimoya"
$ws.Cells.Item(331, 11).Value = "Cultivar V Región"
$ws.Cells.Item(331, 12).Value = "Especial"
$ws.Cells.Item(331, 13).Value = 45
$ws.Cells.Item(331, 14).Value = 31000
$ws.Cells.Item(331, 15).Value = 31000
$ws.Cells.Item(331, 16).Value = 31000
$ws.Cells.Item(331, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(331, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(331, 19).Value = 3100
$ws.Cells.Item(331, 20).Value = 10

# New row 332: Primera, Cultivar V Región, Provincia de Quillota
$ws.Cells.Item(332, 1).Value = 3
$ws.Cells.Item(332, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(332, 3).Value = "Coquimbo"
$ws.Cells.Item(332, 4).Value = 45166
$ws.Cells.Item(332, 5).Value = 5
$ws.Cells.Item(332, 6).Value = "Fruta"
$ws.Cells.Item(332, 7).Value = 100107
$ws.Cells.Item(332, 8).Value = "Otros"
$ws.Cells.Item(332, 9).Value = 100107002
$ws.Cells.Item(332, 10).Value = "Chirimoya"
$ws.Cells.Item(332, 11).Value = "Cultivar V Región"
$ws.Cells.Item(332, 12).Value = "Primera"
$ws.Cells.Item(332, 13).Value = 48
$ws.Cells.Item(332, 14).Value = 29000
$ws.Cells.Item(332, 15).Value = 29000
$ws.Cells.Item(332, 16).Value = 29000
$ws.Cells.Item(332, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(332, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(332, 19).Value = 2900
$ws.Cells.Item(332, 20).Value = 10

# New row 333: Segunda, Cultivar V Región, Provincia de Quillota
$ws.Cells.Item(333, 1).Value = 3
$ws.Cells.Item(333, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(333, 3).Value = "Coquimbo"
$ws.Cells.Item(333, 4).Value = 45166
$ws.Cells.Item(333, 5).Value = 5
$ws.Cells.Item(333, 6).Value = "Fruta"
$ws.Cells.Item(333, 7).Value = 100107
$ws.Cells.Item(333, 8).Value = "Otros"
$ws.Cells.Item(333, 9).Value = 100107002
$ws.Cells.Item(333, 10).Value = "Chirimoya"
$ws.Cells.Item(333, 11).Value = "Cultivar V Región"
$ws.Cells.Item(333, 12).Value = "Segunda"
$ws.Cells.Item(333, 13).Value = 40
$ws.Cells.Item(333, 14).Value = 27000
$ws.Cells.Item(333, 15).Value = 27000
$ws.Cells.Item(333, 16).Value = 27000
$ws.Cells.Item(333, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(333, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(333, 19).Value = 2700
$ws.Cells.Item(333, 20).Value = 10
